$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like plain numbers (e.g. "1.003").
# Force it to Text format before writing so Excel does not coerce these
# inline strings into numeric values, then restore the default style so no
# extra per-cell formatting is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "23.012.52"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.582.30"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "299.09"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").Value = "0.3749"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").Value = "0.3560"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "50.23"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("D10").Value = "1.003"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").Value = "1.208"
$ws.Range("E11").Value = "  -5.60%  "
$ws.Range("D12").Value = "0.07964"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").Value = "21.72"
$ws.Range("E13").Value = "  -6.23%  "
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "7.255"
$ws.Range("E15").Value = "  -5.12%  "
$ws.Range("E16").Value = "  -4.79%  "
$ws.Range("D17").Value = "1.584.99"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "91.65"
$ws.Range("D19").Value = "0.06743"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "6.352"
$ws.Range("E22").Value = "  -3.86%  "
$ws.Range("D23").Value = "23.015.20"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("D25").Value = "2.365"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "2.773"
$ws.Range("E26").Value = "  -5.18%  "
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").Value = "146.59"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").Value = "5.201"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").Value = "131.43"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "2.307"
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").Value = "1.761.36"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").Value = "6.422"
$ws.Range("E33").Value = "  -8.62%  "
$ws.Range("D34").Value = "0.9244"
$ws.Range("E34").Value = "  -6.75%  "
$ws.Range("D35").Value = "0.07306"
$ws.Range("E35").Value = "  -5.81%  "
$ws.Range("D36").Value = "0.02653"
$ws.Range("E36").Value = "  -5.07%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2477"
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "0.08722"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "9.829"
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("D40").Value = "5.946"
$ws.Range("E40").Value = "  -6.20%  "
$ws.Range("D41").Value = "1.328"
$ws.Range("E41").Value = "  -5.19%  "
$ws.Range("D42").Value = "0.6802"
$ws.Range("E42").Value = "  -5.58%  "
$ws.Range("D43").Value = "11.69"
$ws.Range("E43").Value = "  -9.06%  "
$ws.Range("D44").Value = "14.65"
$ws.Range("E44").Value = "  -8.51%  "
$ws.Range("D45").Value = "0.6299"
$ws.Range("E45").Value = "  -5.10%  "
$ws.Range("D46").Value = "3.963"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "2.226"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("D48").Value = "130.28"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").Value = "0.07860"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").Value = "1.173"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "1.156"
$ws.Range("E51").Value = "  -2.78%  "

$priceRange.Style = "Normal"

Write-Output "Applied 97 cell updates"
